$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at AC, shifting old AC onward to the right
$ws.Columns("AC").Insert()

# Set header strings for the new column AC (row 2 = week code, row 3 = date range)
# Order matters for shared string table indices: the date range must be added before the week code
$ws.Range("AC3").Value2 = "06mar-12mar"
$ws.Range("AC2").Value2 = "SE-10"

# Set the data values for the new column
$ws.Range("AC5").Value2 = 1
$ws.Range("AC6").Value2 = 1
$ws.Range("AC8").Value2 = 3
$ws.Range("AC10").Value2 = 1
$ws.Range("AC12").Value2 = 1
$ws.Range("AC14").Value2 = 1

# Fix styles: cells with values should use the "value" style (copy from AB10 which already has that style)
$ws.Range("AB10").Copy()
foreach ($addr in @("AC5","AC6","AC8","AC10","AC12","AC14")) {
    $ws.Range($addr).PasteSpecial(-4122)
}

# Cells without values should use the "blank" style (copy from AB4 which already has that style)
$ws.Range("AB4").Copy()
foreach ($addr in @("AC13","AC15")) {
    $ws.Range($addr).PasteSpecial(-4122)
}

# Set the selection to match
$ws.Range("AD10").Select()
